# Apply the changes described in the commit:
#   1. Switch the table on the "Sources of finance" slide to a different
#      table style ({0E3AE68A-62EA-442F-934B-034A7E0B1D79}).
#   2. Re-colour the deck's theme: the design's colour scheme moves from the
#      "Integral" palette to the "Office Theme" palette (the accompanying
#      theme/fontScheme/fmtScheme stay the same - only the 12 theme colours
#      change).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$newStyleId = "{0E3AE68A-62EA-442F-934B-034A7E0B1D79}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Theme colours -------------------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -> "Office Theme" RGB values
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
